{"js": "// Remove the stale \"_GoBack\" bookmark (it was sitting at the end of the\n// \"Tracking gps soferi\" paragraph) and append three new paragraphs at the\n// end of the document body; Word re-seats \"_GoBack\" around the very last\n// edit location, i.e. a zero-width span at the start of the final\n// paragraph (\"Organigrama\").\n\nconst body = context.document.body;\n\n// 1) Drop the old \"_GoBack\" bookmark wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Find the current last paragraph (\"Programare interviuri\") and append\n//    the three new paragraphs after it, in order.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst p1 = lastParagraph.insertParagraph(\"Evaluari angajati\", \"After\");\nconst p2 = p1.insertParagraph(\"Proceduri interne\", \"After\");\nconst p3 = p2.insertParagraph(\"Organigrama\", \"After\");\nawait context.sync();\n\n// 3) Re-read the paragraph collection so the new last paragraph's anchor is\n//    fresh, then collapse a range at its very start and drop \"_GoBack\"\n//    there (matching Word's own behaviour of re-anchoring the bookmark at\n//    the most recent edit point).\nconst paragraphsAfter = body.paragraphs;\nparagraphsAfter.load(\"items\");\nawait context.sync();\n\nconst newLastParagraph = paragraphsAfter.items[paragraphsAfter.items.length - 1];\nconst startRange = newLastParagraph.getRange(\"Start\");\nstartRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Remove the stale \"_GoBack\" bookmark (it was sitting at the end of the\n# \"Tracking gps soferi\" paragraph) and append three new paragraphs at the\n# end of the document body; Word re-seats \"_GoBack\" around the very last\n# edit location, i.e. a zero-width span at the start of the final\n# paragraph (\"Organigrama\").\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Append the three new paragraphs after the current last paragraph\n#    (\"Programare interviuri\"), in order.\n$r = $d.Paragraphs.Last.Range\n$r.InsertParagraphAfter()\n$r = $d.Paragraphs.Last.Range\n$r.Text = \"Evaluari angajati\"\n\n$r = $d.Paragraphs.Last.Range\n$r.InsertParagraphAfter()\n$r = $d.Paragraphs.Last.Range\n$r.Text = \"Proceduri interne\"\n\n$r = $d.Paragraphs.Last.Range\n$r.InsertParagraphAfter()\n$r = $d.Paragraphs.Last.Range\n$r.Text = \"Organigrama\"\n\n# 3) Re-seat \"_GoBack\" as a collapsed (zero-width) bookmark at the very\n#    start of the new last paragraph (\"Organigrama\"), matching Word's own\n#    behaviour of anchoring the bookmark at the most recent edit point.\n$lastPara = $d.Paragraphs.Last\n$startRange = $lastPara.Range.Duplicate\n$startRange.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $startRange)\n"}
